$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).ColumnWidth = 104.67
$ws.Range("E1").Value = "category"

$st = $wb.Styles.Add("CatStyle")
$st.VerticalAlignment = -4130
$st.HorizontalAlignment = 1
$st.WrapText = $true

$catRange = $ws.Range("E2:E7")
$catRange.Style = "CatStyle"
$catRange.Value = "Профориентация"
Write-Host "done"
